# Auto-generated Excel COM-interop script
# Applies the cell updates described by the diff (F-column 'want to go' counts
# plus two C-column name updates for the STJ01 event).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('展览')
$ws.Range('F3').Value = 1540  # was: 1530
$ws.Range('F4').Value = 805  # was: 803
$ws.Range('F5').Value = 213  # was: 206
$ws.Range('F7').Value = 1115  # was: 1112
$ws.Range('F8').Value = 703  # was: 699
$ws.Range('F9').Value = 761  # was: 760
$ws.Range('F10').Value = 1367  # was: 1358
$ws.Range('F11').Value = 273  # was: 271
$ws.Range('F12').Value = 1016  # was: 1014
$ws.Range('F13').Value = 27  # was: 23
$ws.Range('F14').Value = 59  # was: 57
$ws.Range('F15').Value = 183  # was: 181
$ws.Range('C17').Value = '广州·砂糖桔动漫荟STJ01·综合同人展'  # was: '广州·砂糖桔动漫荟STJ01'
$ws.Range('F17').Value = 427  # was: 418
$ws.Range('F18').Value = 2  # was: 0
$ws.Range('F20').Value = 289  # was: 288
$ws.Range('F21').Value = 535  # was: 534
$ws.Range('F22').Value = 552  # was: 551
$ws.Range('F23').Value = 742  # was: 740
$ws.Range('F24').Value = 227  # was: 224
$ws.Range('F25').Value = 168  # was: 167
$ws.Range('F26').Value = 362  # was: 361

$ws = $wb.Worksheets.Item('演出')
$ws.Range('F3').Value = 987  # was: 983
$ws.Range('F7').Value = 137  # was: 135
$ws.Range('F8').Value = 62  # was: 60
$ws.Range('F9').Value = 582  # was: 581
$ws.Range('F10').Value = 70  # was: 66

$ws = $wb.Worksheets.Item('本地生活')
$ws.Range('F2').Value = 202  # was: 192

$ws = $wb.Worksheets.Item('全部类型')
$ws.Range('F3').Value = 202  # was: 192
$ws.Range('F4').Value = 1540  # was: 1531
$ws.Range('F6').Value = 805  # was: 803
$ws.Range('F7').Value = 213  # was: 206
$ws.Range('F8').Value = 987  # was: 983
$ws.Range('F10').Value = 1115  # was: 1112
$ws.Range('F11').Value = 703  # was: 699
$ws.Range('F12').Value = 761  # was: 760
$ws.Range('F13').Value = 1367  # was: 1358
$ws.Range('F14').Value = 273  # was: 271
$ws.Range('F15').Value = 1016  # was: 1014
$ws.Range('F16').Value = 27  # was: 23
$ws.Range('F17').Value = 59  # was: 0
$ws.Range('F18').Value = 183  # was: 181
$ws.Range('C20').Value = '广州·砂糖桔动漫荟STJ01·综合同人展'  # was: '广州·砂糖桔动漫荟STJ01'
$ws.Range('F20').Value = 427  # was: 418
$ws.Range('F21').Value = 2  # was: 0
$ws.Range('F25').Value = 289  # was: 288
$ws.Range('F27').Value = 137  # was: 135
$ws.Range('F28').Value = 137  # was: 135
$ws.Range('F29').Value = 535  # was: 534
$ws.Range('F30').Value = 552  # was: 551
$ws.Range('F31').Value = 742  # was: 740
$ws.Range('F32').Value = 227  # was: 224
$ws.Range('F33').Value = 62  # was: 60
$ws.Range('F34').Value = 168  # was: 167
$ws.Range('F35').Value = 582  # was: 581
$ws.Range('F36').Value = 70  # was: 66
$ws.Range('F37').Value = 70  # was: 66
$ws.Range('F39').Value = 362  # was: 361
